$d = $word.ActiveDocument

# The 55K shield bolt pattern has now been measured, so remove the two
# outstanding TODO bullets that were tracking that work:
#   "Get accurate dimensions of new 55K shield"
#   "Scrap copper lids thread size"
# Both bullets sit back-to-back right after "Design/make 55K/2.8K stage/shield"
# and right before "Aluminum sheet thickness for rolling our own."

$targets = @(
    "Get accurate dimensions of new 55K shield",
    "Scrap copper lids thread size"
)

foreach ($t in $targets) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $found = $rng.Find.Execute($t, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $para = $rng.Paragraphs(1)
        $para.Range.Delete()
    }
}

# Word keeps the hidden "_GoBack" bookmark pinned to the location of the
# most recent edit. Re-anchor it at the start of the paragraph that now
# immediately follows the deletion ("Aluminum sheet thickness...").
$editRng = $d.Content
$editRng.Find.ClearFormatting()
[void]$editRng.Find.Execute("Aluminum sheet thickness for rolling our own.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$goBackRng = $d.Range($editRng.Start, $editRng.Start)
[void]$d.Bookmarks.Add("_GoBack", $goBackRng)
